# Update "想去人数" (F column) figures across the three sheets that carry
# this data: 展览 (Exhibition), 演出 (Performance) and 全部类型 (All types,
# a merged/filtered view). 本地生活 only has a header row and is untouched.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 10420
$ws1.Range("F3").Value  = 433
$ws1.Range("F4").Value  = 2532
$ws1.Range("F6").Value  = 285
$ws1.Range("F9").Value  = 779
$ws1.Range("F11").Value = 1218
$ws1.Range("F12").Value = 1096
$ws1.Range("F13").Value = 3243
$ws1.Range("F14").Value = 2412
$ws1.Range("F16").Value = 2168
$ws1.Range("F17").Value = 2168
$ws1.Range("F18").Value = 246
$ws1.Range("F19").Value = 1942
$ws1.Range("F22").Value = 578
$ws1.Range("F24").Value = 254
$ws1.Range("F26").Value = 25
$ws1.Range("F27").Value = 245
$ws1.Range("F30").Value = 8
$ws1.Range("F32").Value = 398
$ws1.Range("F34").Value = 20
$ws1.Range("F35").Value = 55
$ws1.Range("F36").Value = 263
$ws1.Range("F37").Value = 9
$ws1.Range("F38").Value = 1580
$ws1.Range("F39").Value = 483
$ws1.Range("F40").Value = 465
$ws1.Range("F41").Value = 1720
$ws1.Range("F42").Value = 141
$ws1.Range("F43").Value = 447
$ws1.Range("F45").Value = 465
$ws1.Range("F46").Value = 1039
$ws1.Range("F48").Value = 367

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 28

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 10420
$ws4.Range("F3").Value  = 433
$ws4.Range("F4").Value  = 2532
$ws4.Range("F5").Value  = 28
$ws4.Range("F8").Value  = 285
$ws4.Range("F11").Value = 779
$ws4.Range("F12").Value = 1097
$ws4.Range("F13").Value = 3243
$ws4.Range("F14").Value = 2412
$ws4.Range("F15").Value = 2168
$ws4.Range("F16").Value = 2168
$ws4.Range("F18").Value = 578
$ws4.Range("F20").Value = 254
$ws4.Range("F22").Value = 25
$ws4.Range("F23").Value = 245
$ws4.Range("F26").Value = 8
$ws4.Range("F28").Value = 398
$ws4.Range("F30").Value = 20
$ws4.Range("F34").Value = 55
$ws4.Range("F35").Value = 263
$ws4.Range("F36").Value = 1580
$ws4.Range("F37").Value = 483
$ws4.Range("F39").Value = 465
$ws4.Range("F40").Value = 1720
$ws4.Range("F41").Value = 141
$ws4.Range("F45").Value = 447
$ws4.Range("F47").Value = 465
$ws4.Range("F48").Value = 1039
$ws4.Range("F49").Value = 367
